$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the B (value) and C (hour) columns for the timeframe-hour MOB auto-update.
$ws.Range("B2").Value = 461.5
$ws.Range("C2").Value = 50

$ws.Range("B3").Value = 251
$ws.Range("C3").Value = 50

$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 50

$ws.Range("B5").Value = 146
$ws.Range("C5").Value = 50

$ws.Range("B6").Value = 30
$ws.Range("C6").Value = 50

$ws.Range("B7").Value = 55
$ws.Range("C7").Value = 50

$ws.Range("B10").Value = 264
$ws.Range("C10").Value = 50

$ws.Range("B11").Value = 150
$ws.Range("C11").Value = 50

$ws.Range("B12").Value = 367.9
$ws.Range("C12").Value = 50

$ws.Range("B13").Value = 344
$ws.Range("C13").Value = 50

$ws.Range("B14").Value = 430
$ws.Range("C14").Value = 50

$ws.Range("B15").Value = 133
$ws.Range("C15").Value = 50

$ws.Range("B17").Value = 75
$ws.Range("C17").Value = 50

$ws.Range("C18").Value = 50

$ws.Range("B19").Value = 6
$ws.Range("C19").Value = 50

$ws.Range("B20").Value = 31
$ws.Range("C20").Value = 50
